$wb = $excel.ActiveWorkbook

# --- Update "Forecast Comparison" sheet ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Cells.Item(2, 2).NumberFormat = "@"
$ws1.Cells.Item(2, 2).Value = "2025-02-02"
$ws1.Cells.Item(2, 2).Style = "Normal"
$ws1.Cells.Item(2, 4).Value = 604
$ws1.Cells.Item(2, 5).Value = 533
$ws1.Cells.Item(2, 6).Value = 572
$ws1.Cells.Item(2, 7).Value = 599
$ws1.Cells.Item(2, 8).Value = 638

$ws1.Cells.Item(3, 2).NumberFormat = "@"
$ws1.Cells.Item(3, 2).Value = "2025-02-09"
$ws1.Cells.Item(3, 2).Style = "Normal"
$ws1.Cells.Item(3, 4).Value = 638
$ws1.Cells.Item(3, 5).Value = 541
$ws1.Cells.Item(3, 6).Value = 591
$ws1.Cells.Item(3, 7).Value = 628
$ws1.Cells.Item(3, 8).Value = 681

$ws1.Cells.Item(4, 2).NumberFormat = "@"
$ws1.Cells.Item(4, 2).Value = "2025-02-16"
$ws1.Cells.Item(4, 2).Style = "Normal"
$ws1.Cells.Item(4, 4).Value = 667
$ws1.Cells.Item(4, 5).Value = 527
$ws1.Cells.Item(4, 6).Value = 581
$ws1.Cells.Item(4, 7).Value = 622
$ws1.Cells.Item(4, 8).Value = 680

$ws1.Cells.Item(5, 2).NumberFormat = "@"
$ws1.Cells.Item(5, 2).Value = "2025-02-23"
$ws1.Cells.Item(5, 2).Style = "Normal"
$ws1.Cells.Item(5, 4).Value = 672
$ws1.Cells.Item(5, 5).Value = 525
$ws1.Cells.Item(5, 6).Value = 580
$ws1.Cells.Item(5, 7).Value = 621
$ws1.Cells.Item(5, 8).Value = 681

$ws1.Cells.Item(6, 2).NumberFormat = "@"
$ws1.Cells.Item(6, 2).Value = "2025-03-02"
$ws1.Cells.Item(6, 2).Style = "Normal"
$ws1.Cells.Item(6, 4).Value = 652
$ws1.Cells.Item(6, 5).Value = 530
$ws1.Cells.Item(6, 6).Value = 588
$ws1.Cells.Item(6, 7).Value = 632
$ws1.Cells.Item(6, 8).Value = 696

$ws1.Cells.Item(7, 2).NumberFormat = "@"
$ws1.Cells.Item(7, 2).Value = "2025-03-09"
$ws1.Cells.Item(7, 2).Style = "Normal"
$ws1.Cells.Item(7, 4).Value = 625
$ws1.Cells.Item(7, 5).Value = 527
$ws1.Cells.Item(7, 6).Value = 588
$ws1.Cells.Item(7, 7).Value = 635
$ws1.Cells.Item(7, 8).Value = 703

$ws1.Cells.Item(8, 2).NumberFormat = "@"
$ws1.Cells.Item(8, 2).Value = "2025-03-16"
$ws1.Cells.Item(8, 2).Style = "Normal"
$ws1.Cells.Item(8, 4).Value = 610
$ws1.Cells.Item(8, 5).Value = 501
$ws1.Cells.Item(8, 6).Value = 562
$ws1.Cells.Item(8, 7).Value = 610
$ws1.Cells.Item(8, 8).Value = 681

$ws1.Cells.Item(9, 2).NumberFormat = "@"
$ws1.Cells.Item(9, 2).Value = "2025-03-23"
$ws1.Cells.Item(9, 2).Style = "Normal"
$ws1.Cells.Item(9, 4).Value = 619
$ws1.Cells.Item(9, 5).Value = 528
$ws1.Cells.Item(9, 6).Value = 590
$ws1.Cells.Item(9, 7).Value = 638
$ws1.Cells.Item(9, 8).Value = 709

$ws1.Cells.Item(10, 2).NumberFormat = "@"
$ws1.Cells.Item(10, 2).Value = "2025-03-30"
$ws1.Cells.Item(10, 2).Style = "Normal"
$ws1.Cells.Item(10, 4).Value = 643
$ws1.Cells.Item(10, 5).Value = 505
$ws1.Cells.Item(10, 6).Value = 570
$ws1.Cells.Item(10, 7).Value = 621
$ws1.Cells.Item(10, 8).Value = 697

$ws1.Cells.Item(11, 2).NumberFormat = "@"
$ws1.Cells.Item(11, 2).Value = "2025-04-06"
$ws1.Cells.Item(11, 2).Style = "Normal"
$ws1.Cells.Item(11, 4).Value = 655
$ws1.Cells.Item(11, 5).Value = 504
$ws1.Cells.Item(11, 6).Value = 569
$ws1.Cells.Item(11, 7).Value = 620
$ws1.Cells.Item(11, 8).Value = 695

$ws1.Cells.Item(12, 2).NumberFormat = "@"
$ws1.Cells.Item(12, 2).Value = "2025-04-13"
$ws1.Cells.Item(12, 2).Style = "Normal"
$ws1.Cells.Item(12, 4).Value = 660
$ws1.Cells.Item(12, 5).Value = 508
$ws1.Cells.Item(12, 6).Value = 572
$ws1.Cells.Item(12, 7).Value = 622
$ws1.Cells.Item(12, 8).Value = 695

$ws1.Cells.Item(13, 2).NumberFormat = "@"
$ws1.Cells.Item(13, 2).Value = "2025-04-20"
$ws1.Cells.Item(13, 2).Style = "Normal"
$ws1.Cells.Item(13, 4).Value = 647
$ws1.Cells.Item(13, 5).Value = 498
$ws1.Cells.Item(13, 6).Value = 562
$ws1.Cells.Item(13, 7).Value = 613
$ws1.Cells.Item(13, 8).Value = 689

$ws1.Cells.Item(14, 2).NumberFormat = "@"
$ws1.Cells.Item(14, 2).Value = "2025-04-27"
$ws1.Cells.Item(14, 2).Style = "Normal"
$ws1.Cells.Item(14, 4).Value = 641
$ws1.Cells.Item(14, 5).Value = 493
$ws1.Cells.Item(14, 6).Value = 558
$ws1.Cells.Item(14, 7).Value = 610
$ws1.Cells.Item(14, 8).Value = 686

$ws1.Cells.Item(15, 2).NumberFormat = "@"
$ws1.Cells.Item(15, 2).Value = "2025-05-04"
$ws1.Cells.Item(15, 2).Style = "Normal"
$ws1.Cells.Item(15, 4).Value = 623
$ws1.Cells.Item(15, 5).Value = 479
$ws1.Cells.Item(15, 6).Value = 544
$ws1.Cells.Item(15, 7).Value = 597
$ws1.Cells.Item(15, 8).Value = 675

$ws1.Cells.Item(16, 2).NumberFormat = "@"
$ws1.Cells.Item(16, 2).Value = "2025-05-11"
$ws1.Cells.Item(16, 2).Style = "Normal"
$ws1.Cells.Item(16, 4).Value = 627
$ws1.Cells.Item(16, 5).Value = 482
$ws1.Cells.Item(16, 6).Value = 547
$ws1.Cells.Item(16, 7).Value = 599
$ws1.Cells.Item(16, 8).Value = 676

$ws1.Cells.Item(17, 2).NumberFormat = "@"
$ws1.Cells.Item(17, 2).Value = "2025-05-18"
$ws1.Cells.Item(17, 2).Style = "Normal"
$ws1.Cells.Item(17, 4).Value = 623
$ws1.Cells.Item(17, 5).Value = 479
$ws1.Cells.Item(17, 6).Value = 545
$ws1.Cells.Item(17, 7).Value = 599
$ws1.Cells.Item(17, 8).Value = 679

# --- Update "Summary" sheet ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Cells.Item(2, 2).NumberFormat = "@"
$ws2.Cells.Item(2, 2).Value = "2022-12-25 to 2025-01-26"
$ws2.Cells.Item(2, 2).Style = "Normal"

$ws2.Cells.Item(4, 2).NumberFormat = "@"
$ws2.Cells.Item(4, 2).Value = "1124"
$ws2.Cells.Item(4, 2).Style = "Normal"

$ws2.Cells.Item(6, 2).NumberFormat = "@"
$ws2.Cells.Item(6, 2).Value = "560"
$ws2.Cells.Item(6, 2).Style = "Normal"

$ws2.Cells.Item(7, 2).NumberFormat = "@"
$ws2.Cells.Item(7, 2).Value = "230"
$ws2.Cells.Item(7, 2).Style = "Normal"

$ws2.Cells.Item(8, 2).NumberFormat = "@"
$ws2.Cells.Item(8, 2).Value = "56512 units"
$ws2.Cells.Item(8, 2).Style = "Normal"

$ws2.Cells.Item(9, 2).NumberFormat = "@"
$ws2.Cells.Item(9, 2).Value = "10206"
$ws2.Cells.Item(9, 2).Style = "Normal"

$ws2.Cells.Item(10, 2).NumberFormat = "@"
$ws2.Cells.Item(10, 2).Value = "5087"
$ws2.Cells.Item(10, 2).Style = "Normal"

$ws2.Cells.Item(11, 2).NumberFormat = "@"
$ws2.Cells.Item(11, 2).Value = "2581"
$ws2.Cells.Item(11, 2).Style = "Normal"

$ws2.Cells.Item(12, 2).NumberFormat = "@"
$ws2.Cells.Item(12, 2).Value = "672"
$ws2.Cells.Item(12, 2).Style = "Normal"

$ws2.Cells.Item(13, 2).NumberFormat = "@"
$ws2.Cells.Item(13, 2).Value = "2025-02-23"
$ws2.Cells.Item(13, 2).Style = "Normal"

$ws2.Cells.Item(14, 2).NumberFormat = "@"
$ws2.Cells.Item(14, 2).Value = "604"
$ws2.Cells.Item(14, 2).Style = "Normal"

$ws2.Cells.Item(15, 2).NumberFormat = "@"
$ws2.Cells.Item(15, 2).Value = "2025-02-02"
$ws2.Cells.Item(15, 2).Style = "Normal"
